$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns (row 1) from descriptive Spanish labels to short field codes
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Normalize capitalization of state/municipality names: title-case every
#    word (e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga"), matching
#    Excel's PROPER() behavior, for all data rows in columns A and B.
for ($r = 2; $r -le 1241; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $txt = $cell.Text
        if ($txt -ne "") {
            $cell.Value = $excel.WorksheetFunction.Proper($txt)
        }
    }
}

# 3. Remove the trailing source/footnote rows (1243-1247), shrinking the
#    sheet's used range back down to row 1241.
$ws.Range("A1243:D1247").Delete(-4162) | Out-Null
